$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.MoveEnd(1, -1)
    $r.Text = $newText
}

function Set-CellBold($table, $row, $col, $boldValue) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $start = $r.Start
    $len = ($r.End - 1) - $r.Start
    $r2 = $d.Range($start, $start + $len)
    $r2.Font.Bold = $boldValue
}

# Row 4 = "Neg-Bin. GAM"
Set-CellText $t 4 3 "142"
Set-CellText $t 4 4 "-28.3"
Set-CellText $t 4 5 "2.6"
Set-CellText $t 4 6 "120"
Set-CellText $t 4 7 "-0.2"
Set-CellBold $t 4 7 1

# Row 5 = "Poisson GAM"
Set-CellText $t 5 4 "-31.9"
Set-CellText $t 5 5 "2.9"
Set-CellText $t 5 6 "134"
Set-CellText $t 5 7 "-3.3"

# Row 6 = "Serfling-Poisson GLM"
Set-CellText $t 6 3 "175"
Set-CellText $t 6 4 "-26.4"
Set-CellText $t 6 6 "141"
Set-CellText $t 6 7 "-1.7"
Set-CellBold $t 6 7 0
